# Commit: "MOD:  - Neigungssensor TestIBN  - Messdaten Fahrwer auf Distanz in Excel"
#
#   - Add a new worksheet "Fahrversuche" (after "Tabelle1") holding the
#     drive-test distance measurements, and make it the active tab.
#   - Tabelle1 keeps all of its data/formulas; only its view state changes
#     (it is no longer the selected tab, and its scroll position moves).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tabelle1")

# --- New sheet with the "Fahrversuche" measurement data, placed after Tabelle1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Fahrversuche"

# Header row
$ws2.Range("A1").Value = "delay"
$ws2.Range("B1").Value = "Speed"
$ws2.Range("C1").Value = "Richtung"
$ws2.Range("D1").Value = "Distanz gefahren [cm]"

# Row 2 - baseline delay/speed plus first direction's measurements
$ws2.Range("A2").Value = 4100
$ws2.Range("B2").Value = "SPEED_GANZLANGSAM"
$ws2.Range("C2").Value = "geradeaus hoch"
$ws2.Range("D2").Value = 7.5
$ws2.Range("E2").Value = 8.5

# Rows 3-9 - one direction + its measured runs per row
$ws2.Range("C3").Value = "geradeaus runter"
$ws2.Range("D3").Value = 11
$ws2.Range("E3").Value = 11

$ws2.Range("C4").Value = "links nach rechts"
$ws2.Range("D4").Value = 10
$ws2.Range("E4").Value = 9
$ws2.Range("F4").Value = 9.5

$ws2.Range("C5").Value = "rechts nach links"
$ws2.Range("D5").Value = 11.5
$ws2.Range("E5").Value = 11.5

$ws2.Range("C6").Value = "diagonal links oben rechts unten"
$ws2.Range("D6").Value = 10
$ws2.Range("E6").Value = 10

$ws2.Range("C7").Value = "diagonal rechts oben links unten"
$ws2.Range("D7").Value = 9.5
$ws2.Range("E7").Value = 9.8000000000000007

$ws2.Range("C8").Value = "diagonal links unten rechts oben"
$ws2.Range("D8").Value = 11
$ws2.Range("E8").Value = 11.5

$ws2.Range("C9").Value = "diagonal rechts unten links oben"
$ws2.Range("D9").Value = 10.5
$ws2.Range("E9").Value = 11

# Column widths so the header/direction text fits (B & C were auto-fit in the
# original; D was widened manually). ColumnWidth is quantized by the host to
# the nearest 1/6th of a character, so these inputs are chosen to land on the
# closest achievable width to the authored 21.5703125 / 30.140625 / 30.85546875.
$ws2.Columns.Item(2).ColumnWidth = 20.666666666666668
$ws2.Columns.Item(3).ColumnWidth = 29.333333333333332
$ws2.Columns.Item(4).ColumnWidth = 30

# Selection left on the new sheet
$ws2.Range("C8").Select() | Out-Null

# --- Tabelle1's view scrolls up a bit, selection stays on C30 ---
$ws1.Range("A8").Select() | Out-Null
$ws1.Range("C30").Select() | Out-Null

# Fahrversuche ends up the active/selected tab
$ws2.Activate() | Out-Null
